$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of mod-count data appended below the existing table (row 14 -> row 15)
$newRow = 15

# Write the raw values first. The date column is plain text in this sheet
# (e.g. "2025/11/23"), so force Text formatting before assignment to stop
# Excel's automatic date-literal recognition from turning it into a serial
# number.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025/11/24"
$ws.Range("B" + $newRow).Value = "逃离鸭科夫"
$ws.Range("C" + $newRow).Value = 1242

# Match the formatting of the previous data row (centered alignment, same
# number format) by copying its style onto the new row.
$ws.Range("A14:C14").Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)
